# Final versions submitted to the journal now
#
# The reviewer-response document still carried two tracked "formatting
# changed" revisions (rPrChange) recorded by Andy Wills back when the
# bold emphasis was added to his two inline replies. Now that this is
# the final version being submitted, those tracked changes are accepted
# (the new bold/italic/colour formatting is kept, the rPrChange history
# of what the formatting used to be is discarded).

$d = $word.ActiveDocument

# Walk the revisions back-to-front and accept every one of them. Using
# the individual Revision objects (rather than Revisions.AcceptAll /
# Document.AcceptAllRevisions) keeps everything else about each run
# untouched -- only the w:rPrChange markup that records the "previous"
# run formatting is removed, exactly as if the author clicked "Accept"
# on each tracked change in turn.
$revisions = $d.Revisions
for ($i = $revisions.Count; $i -ge 1; $i--) {
    $revisions.Item($i).Accept()
}

# No outstanding comments in this document, but accept/clear them too
# for good measure -- once a document is being finalised for submission
# there should be no reviewer markup of any kind left in it.
$comments = $d.Comments
for ($i = $comments.Count; $i -ge 1; $i--) {
    $comments.Item($i).Delete()
}

# Stop recording further changes now that the document is final.
$d.TrackRevisions = $false
